$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("C8").Value = 16
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4"
$ws.Range("E8").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F8").Value = 662
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "10592.00"

# --- Row 9 ---
$ws.Range("A9").Value = ""
$ws.Range("C9").Value = 75
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.0"
$ws.Range("E9").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.00"

# --- Row 10 ---
$ws.Range("A10").Value = "Mtr."
$ws.Range("C10").Value = 71
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19"
$ws.Range("E10").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F10").Value = 81
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "5751.00"

# --- Row 11 ---
$ws.Range("C11").Value = 96
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31"
$ws.Range("E11").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'

# --- Row 12 ---
$ws.Range("C12").Value = 10
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32"
$ws.Range("E12").Value = ' 50/63 A rating'
$ws.Range("F12").Value = 900
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "9000.00"

# --- Row 14 (Grand Total) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "25343.00"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "25343.00"

# --- Row 16 (Net Payable Amount) ---
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "25343.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "25343.00"
